$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row data (after the edit): Column A = Tag, Column B = Code, Column C = Type (always 1)
# Order below matches the row order top-to-bottom as they appear after the edit (rows 2..35).
$rows = @(
    @("haircut", "G001"),
    @("general hair cut", "G001"),
    @("cut", "G001"),
    @("children haircut", "G001-1"),
    @("childrens haircut", "G001-1"),
    @("children's haircut", "G001-1"),
    @("men's haircut", "G001-2"),
    @("mens haircut", "G001-2"),
    @("men haircut", "G001-2"),
    @("women's haircut", "G001-3"),
    @("womens haircut", "G001-3"),
    @("women haircut", "G001-3"),
    @("ladies haircut", "G001-3"),
    @("ladies haircut", "G001-3"),
    @("hair setup", "G002"),
    @("hair setups", "G002"),
    @("setup", "G002"),
    @("setups", "G002"),
    @("hair dressings", "G003"),
    @("hair dressing", "G003"),
    @("dressings", "G003"),
    @("dressing", "G003"),
    @("massage", "G004"),
    @("head massage", "G004"),
    @("scalp massage", "G004"),
    @("hair massage", "G004"),
    @("facial", "G005"),
    @("facials", "G005"),
    @("pedicure", "G006"),
    @("manicure", "G007"),
    @("hair colouring", "G008"),
    @("hair coloring", "G008"),
    @("colouring", "G008"),
    @("coloring", "G008")
)

# Clear out the previously existing data rows below the header before rewriting,
# since the new table is longer than the old one (old: rows 2-22, new: rows 2-35).
$ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(35, 3)).Clear()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = 1
    $r = $r + 1
}

# Update selection / active cell to match the final saved view state.
$ws.Range("F14").Select()

$wb.Save()
